$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D31").Value = 9938.521000000001
$ws.Range("D32").Value = 8865.171
$ws.Range("D33").Value = 8453.887000000001
$ws.Range("D34").Value = 14397.439
$ws.Range("D35").Value = 10743.533
$ws.Range("D36").Value = 10126.608
$ws.Range("D37").Value = 9489.620000000001
$ws.Range("D38").Value = 15004.331
